# Auto-generated edit script: apply Gungnir_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10936
$ws.Range("H23").Value = 10000
$ws.Range("J23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10468
$ws.Range("H38").Value = 1035.64
$ws.Range("I38").Value = 98.2
$ws.Range("J38").Value = 1660.6
$ws.Range("K38").Value = 294.6
$ws.Range("L38").Value = 4981.799999999999
$ws.Range("M38").Value = 77.39999999999998
$ws.Range("N38").Value = -5725.799999999999
$ws.Range("H43").Value = 27779786
$ws.Range("J43").Value = 1734.5454
$ws.Range("L43").Value = 1734.5454
$ws.Range("N43").Value = -1872.5454
$ws.Range("H64").Value = 3332.889
$ws.Range("I64").Value = 2999
$ws.Range("K64").Value = 2999
$ws.Range("M64").Value = -2751
$ws.Range("H67").Value = 3332.889
$ws.Range("I67").Value = 2999
$ws.Range("K67").Value = 2999
$ws.Range("M67").Value = -2141
$ws.Range("H69").Value = 1000
$ws.Range("I69").Value = 1000
$ws.Range("K69").Value = 3000
$ws.Range("M69").Value = -2126
$ws.Range("H72").Value = 1000
$ws.Range("I72").Value = 1000
$ws.Range("K72").Value = 9000
$ws.Range("M72").Value = -4632
$ws.Range("H121").Value = 628.1786
$ws.Range("J121").Value = 677.0833
$ws.Range("L121").Value = 2031.2499
$ws.Range("N121").Value = -5525.2499
$ws.Range("H125").Value = 965.7368
$ws.Range("I125").Value = 831.3333
$ws.Range("J125").Value = 1086.7
$ws.Range("K125").Value = 7481.9997
$ws.Range("L125").Value = 9780.300000000001
$ws.Range("M125").Value = -5021.9997
$ws.Range("N125").Value = -14700.3
$ws.Range("H131").Value = 3658
$ws.Range("I131").Value = 747.5
$ws.Range("J131").Value = 5598.3335
$ws.Range("K131").Value = 2242.5
$ws.Range("L131").Value = 16795.0005
$ws.Range("M131").Value = 2797.5
$ws.Range("N131").Value = -26875.0005
$ws.Range("H138").Value = 3509.1
$ws.Range("I138").Value = 1856.1154
$ws.Range("J138").Value = 4089.8784
$ws.Range("K138").Value = 5568.3462
$ws.Range("L138").Value = 12269.6352
$ws.Range("M138").Value = -428.3462
$ws.Range("N138").Value = -22549.6352
$ws.Range("H141").Value = 5259.222
$ws.Range("I141").Value = 2297.0667
$ws.Range("J141").Value = 20070
$ws.Range("K141").Value = 6891.2001
$ws.Range("L141").Value = 60210
$ws.Range("M141").Value = -1711.2001
$ws.Range("N141").Value = -70570

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 115.2
$ws.Range("I5").Value = 81.57143000000001
$ws.Range("J5").Value = 193.66667
$ws.Range("K5").Value = 81.57143000000001
$ws.Range("L5").Value = 193.66667
$ws.Range("M5").Value = 30.42856999999999
$ws.Range("N5").Value = -417.66667
$ws.Range("H61").Value = 1848.9574
$ws.Range("I61").Value = 1818.7812
$ws.Range("J61").Value = 1913.3334
$ws.Range("K61").Value = 1818.7812
$ws.Range("L61").Value = 1913.3334
$ws.Range("M61").Value = -1606.7812
$ws.Range("N61").Value = -2337.3334
$ws.Range("H74").Value = 1630.0358
$ws.Range("I74").Value = 1773.5385
$ws.Range("K74").Value = 1773.5385
$ws.Range("M74").Value = -899.5385000000001
$ws.Range("H77").Value = 1630.0358
$ws.Range("I77").Value = 1773.5385
$ws.Range("K77").Value = 8867.692500000001
$ws.Range("M77").Value = -4499.692500000001
$ws.Range("H136").Value = 1848.9574
$ws.Range("I136").Value = 1818.7812
$ws.Range("J136").Value = 1913.3334
$ws.Range("K136").Value = 5456.3436
$ws.Range("L136").Value = 5740.0002
$ws.Range("M136").Value = -2906.3436
$ws.Range("N136").Value = -10840.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 115.2
$ws.Range("I4").Value = 81.57143000000001
$ws.Range("J4").Value = 193.66667
$ws.Range("K4").Value = 81.57143000000001
$ws.Range("L4").Value = 193.66667
$ws.Range("M4").Value = 33.42856999999999
$ws.Range("N4").Value = -423.66667
$ws.Range("H62").Value = 19900
$ws.Range("J62").Value = 19900
$ws.Range("L62").Value = 19900
$ws.Range("N62").Value = -21272
$ws.Range("H65").Value = 19900
$ws.Range("J65").Value = 19900
$ws.Range("L65").Value = 59700
$ws.Range("N65").Value = -66564
$ws.Range("H105").Value = 25001936
$ws.Range("I105").Value = 1897.5143
$ws.Range("J105").Value = 200002200
$ws.Range("K105").Value = 1897.5143
$ws.Range("L105").Value = 200002200
$ws.Range("M105").Value = -150.5143
$ws.Range("N105").Value = -200005694

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 22954.637
$ws.Range("I7").Value = 83.5
$ws.Range("J7").Value = 50400
$ws.Range("K7").Value = 83.5
$ws.Range("L7").Value = 50400
$ws.Range("M7").Value = 29.5
$ws.Range("N7").Value = -50626
$ws.Range("H134").Value = 1691.2122
$ws.Range("I134").Value = 1719.6774
$ws.Range("J134").Value = 1250
$ws.Range("K134").Value = 5159.0322
$ws.Range("L134").Value = 3750
$ws.Range("M134").Value = -2624.0322
$ws.Range("N134").Value = -8820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 5458.3335
$ws.Range("I33").Value = 1168.6666
$ws.Range("J33").Value = 7174.2
$ws.Range("K33").Value = 7011.9996
$ws.Range("L33").Value = 43045.2
$ws.Range("M33").Value = -6728.9996
$ws.Range("N33").Value = -43611.2
$ws.Range("H44").Value = 83335496
$ws.Range("I44").Value = 125002750
$ws.Range("J44").Value = 992
$ws.Range("K44").Value = 375008250
$ws.Range("L44").Value = 2976
$ws.Range("M44").Value = -375007852
$ws.Range("N44").Value = -3772
$ws.Range("H64").Value = 1298
$ws.Range("I64").Value = 953
$ws.Range("J64").Value = 1360.7273
$ws.Range("K64").Value = 2859
$ws.Range("L64").Value = 4082.1819
$ws.Range("M64").Value = -2589
$ws.Range("N64").Value = -4622.1819
$ws.Range("H67").Value = 1298
$ws.Range("I67").Value = 953
$ws.Range("J67").Value = 1360.7273
$ws.Range("K67").Value = 2859
$ws.Range("L67").Value = 4082.1819
$ws.Range("M67").Value = -1923
$ws.Range("N67").Value = -5954.1819
$ws.Range("H113").Value = 11458930
$ws.Range("I113").Value = 10417167
$ws.Range("J113").Value = 12500693
$ws.Range("K113").Value = 31251501
$ws.Range("L113").Value = 37502079
$ws.Range("M113").Value = -31249331
$ws.Range("N113").Value = -37506419

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 100000140
$ws.Range("I2").Value = 142857180
$ws.Range("J2").Value = 406.66666
$ws.Range("K2").Value = 142857180
$ws.Range("L2").Value = 406.66666
$ws.Range("M2").Value = -142857067
$ws.Range("N2").Value = -632.66666
$ws.Range("H126").Value = 4129.0586
$ws.Range("I126").Value = 3171.2856
$ws.Range("J126").Value = 4799.5
$ws.Range("K126").Value = 9513.856800000001
$ws.Range("L126").Value = 14398.5
$ws.Range("M126").Value = -7043.856800000001
$ws.Range("N126").Value = -19338.5
$ws.Range("H131").Value = 24730.4
$ws.Range("J131").Value = 24730.4
$ws.Range("L131").Value = 24730.4
$ws.Range("N131").Value = -34810.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3250
$ws.Range("I7").Value = 3500
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 3500
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -3388
$ws.Range("N7").Value = -3224
$ws.Range("H40").Value = 250000000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 250000000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 250000000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -250000272
$ws.Range("H46").Value = 3472727.5
$ws.Range("I46").Value = 5208752
$ws.Range("J46").Value = 678.5
$ws.Range("K46").Value = 5208752
$ws.Range("L46").Value = 678.5
$ws.Range("M46").Value = -5208564
$ws.Range("N46").Value = -1054.5
$ws.Range("H74").Value = 29444.5
$ws.Range("J74").Value = 25926
$ws.Range("L74").Value = 25926
$ws.Range("N74").Value = -27922
$ws.Range("H77").Value = 29444.5
$ws.Range("J77").Value = 25926
$ws.Range("L77").Value = 77778
$ws.Range("N77").Value = -87762
$ws.Range("H82").Value = 1609.8
$ws.Range("I82").Value = 1574.75
$ws.Range("J82").Value = 1633.1666
$ws.Range("K82").Value = 1574.75
$ws.Range("L82").Value = 1633.1666
$ws.Range("M82").Value = -1213.75
$ws.Range("N82").Value = -2355.1666
$ws.Range("H85").Value = 1609.8
$ws.Range("I85").Value = 1574.75
$ws.Range("J85").Value = 1633.1666
$ws.Range("K85").Value = 1574.75
$ws.Range("L85").Value = 1633.1666
$ws.Range("M85").Value = -326.75
$ws.Range("N85").Value = -4129.1666
$ws.Range("H122").Value = 52000
$ws.Range("I122").Value = 62500
$ws.Range("K122").Value = 187500
$ws.Range("M122").Value = -185050
$ws.Range("H126").Value = 3250
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -13940
$ws.Range("H136").Value = 6096.4
$ws.Range("I136").Value = 4466.0586
$ws.Range("K136").Value = 13398.1758
$ws.Range("M136").Value = -10848.1758

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9800
$ws.Range("J15").Value = 9800
$ws.Range("L15").Value = 9800
$ws.Range("N15").Value = -10376

Write-Output "edits applied"